$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "62.921.60"
$ws.Range("E2").Value2 = "  +2.26%  "

$ws.Range("D3").Value2 = "3.478.71"
$ws.Range("E3").Value2 = "  +2.57%  "

$ws.Range("E4").Value2 = "  -0.02%  "

$cell = $ws.Range("D5")
$cell.Value2 = "'582.18"
$cell.ClearFormats()
$ws.Range("E5").Value2 = "  +1.15%  "

$cell = $ws.Range("D6")
$cell.Value2 = "'147.65"
$cell.ClearFormats()
$ws.Range("E6").Value2 = "  +4.69%  "

$ws.Range("E7").Value2 = "  -0.10%  "

$ws.Range("E8").Value2 = "  +1.43%  "

$cell = $ws.Range("D9")
$cell.Value2 = "'7.67"
$cell.ClearFormats()
$ws.Range("E9").Value2 = "  -0.85%  "

$cell = $ws.Range("D10")
$cell.Value2 = "'0.125"
$cell.ClearFormats()
$ws.Range("E10").Value2 = "  +2.33%  "

$ws.Range("E11").Value2 = "  +3.54%  "

$ws.Range("D12").Value2 = "4.075.95"
$ws.Range("E12").Value2 = "  +2.60%  "

$cell = $ws.Range("D13")
$cell.Value2 = "'29.95"
$cell.ClearFormats()
$ws.Range("E13").Value2 = "  +5.45%  "

$ws.Range("E14").Value2 = "  +0.57%  "

$ws.Range("D15").Value2 = "3.461.92"
$ws.Range("E15").Value2 = "  +2.15%  "

$ws.Range("E16").Value2 = "  +1.02%  "

$ws.Range("D17").Value2 = "62.947.12"
$ws.Range("E17").Value2 = "  +2.23%  "

$cell = $ws.Range("D18")
$cell.Value2 = "'6.38"
$cell.ClearFormats()
$ws.Range("E18").Value2 = "  +3.80%  "

$cell = $ws.Range("D19")
$cell.Value2 = "'14.41"
$cell.ClearFormats()
$ws.Range("E19").Value2 = "  +5.40%  "

$ws.Range("E20").Value2 = "  +4.08%  "

$cell = $ws.Range("D21")
$cell.Value2 = "'390.42"
$cell.ClearFormats()
$ws.Range("E21").Value2 = "  -0.15%  "

$cell = $ws.Range("D22")
$cell.Value2 = "'0.568"
$cell.ClearFormats()
$ws.Range("E22").Value2 = "  +2.58%  "

$cell = $ws.Range("D23")
$cell.Value2 = "'75.19"
$cell.ClearFormats()
$ws.Range("E23").Value2 = "  -0.19%  "

$ws.Range("E24").Value2 = "  -0.03%  "

$ws.Range("D25").Value2 = "3.624.55"
$ws.Range("E25").Value2 = "  +2.57%  "

$ws.Range("E26").Value2 = "  +3.29%  "

$cell = $ws.Range("D27")
$cell.Value2 = "'0.182"
$cell.ClearFormats()
$ws.Range("E27").Value2 = "  -5.50%  "

$cell = $ws.Range("D28")
$cell.Value2 = "'7.69"
$cell.ClearFormats()
$ws.Range("E28").Value2 = "  +5.90%  "

$cell = $ws.Range("D29")
$cell.Value2 = "'0.999"
$cell.ClearFormats()
$ws.Range("E29").Value2 = "  +0.02%  "

$ws.Range("E30").Value2 = "  +2.68%  "

$cell = $ws.Range("D31")
$cell.Value2 = "'2.15"
$cell.ClearFormats()
$ws.Range("E31").Value2 = "  +0.13%  "

$ws.Range("E32").Value2 = "  +2.24%  "

$ws.Range("E33").Value2 = "  +0.02%  "

$cell = $ws.Range("D34")
$cell.Value2 = "'23.87"
$cell.ClearFormats()
$ws.Range("E34").Value2 = "  +2.44%  "

$ws.Range("E35").Value2 = "  +2.93%  "

$cell = $ws.Range("D36")
$cell.Value2 = "'5.30"
$cell.ClearFormats()
$ws.Range("E36").Value2 = "  +4.91%  "

$cell = $ws.Range("D37")
$cell.Value2 = "'31.70"
$cell.ClearFormats()
$ws.Range("E37").Value2 = "  +21.94%  "

$cell = $ws.Range("D38")
$cell.Value2 = "'171.51"
$cell.ClearFormats()
$ws.Range("E38").Value2 = "  +2.06%  "

$cell = $ws.Range("D39")
$cell.Value2 = "'1.57"
$cell.ClearFormats()
$ws.Range("E39").Value2 = "  +7.03%  "

$ws.Range("D40").Value2 = "3.518.00"
$ws.Range("E40").Value2 = "  +2.64%  "

$cell = $ws.Range("D41")
$cell.Value2 = "'0.0771"
$cell.ClearFormats()
$ws.Range("E41").Value2 = "  -0.15%  "

$ws.Range("E42").Value2 = "  +3.39%  "

$cell = $ws.Range("D43")
$cell.Value2 = "'42.29"
$cell.ClearFormats()
$ws.Range("E43").Value2 = "  -0.47%  "

$cell = $ws.Range("D44")
$cell.Value2 = "'4.47"
$cell.ClearFormats()
$ws.Range("E44").Value2 = "  +1.03%  "

$ws.Range("E45").Value2 = "  +3.58%  "

$ws.Range("E46").Value2 = "  +4.60%  "

$ws.Range("D47").Value2 = "2.607.43"
$ws.Range("E47").Value2 = "  +5.34%  "

$cell = $ws.Range("D48")
$cell.Value2 = "'23.60"
$cell.ClearFormats()
$ws.Range("E48").Value2 = "  +2.56%  "

$cell = $ws.Range("D49")
$cell.Value2 = "'2.27"
$cell.ClearFormats()
$ws.Range("E49").Value2 = "  +10.05%  "

$cell = $ws.Range("D50")
$cell.Value2 = "'6.80"
$cell.ClearFormats()
$ws.Range("E50").Value2 = "  +1.99%  "

$ws.Range("B51").Value2 = "VeChain"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D51")
$cell.Value2 = "'0.0269"
$cell.ClearFormats()
$ws.Range("E51").Value2 = "  +2.46%  "
